$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.501.90'
$ws.Range("E2").Value = '  -2.18%  '
$ws.Range("D3").Value = '1.749.46'
$ws.Range("E3").Value = '  -2.32%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.51'
$ws.Range("E5").Value = '  +0.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4456'
$ws.Range("E7").Value = '  +3.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3602'
$ws.Range("E8").Value = '  -0.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07483'
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.88'
$ws.Range("E10").Value = '  -6.33%  '
$ws.Range("E11").Value = '  -1.93%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.59'
$ws.Range("E13").Value = '  -4.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.018'
$ws.Range("E14").Value = '  -2.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.135'
$ws.Range("E15").Value = '  -2.48%  '
$ws.Range("D16").Value = '1.755.50'
$ws.Range("E16").Value = '  -2.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.67'
$ws.Range("E17").Value = '  +0.46%  '
$ws.Range("E18").Value = '  -0.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06409'
$ws.Range("E19").Value = '  +0.88%  '
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.80'
$ws.Range("E21").Value = '  -2.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.843'
$ws.Range("E22").Value = '  -2.36%  '
$ws.Range("D23").Value = '27.548.50'
$ws.Range("E23").Value = '  -2.02%  '
$ws.Range("E24").Value = '  -2.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.101'
$ws.Range("E25").Value = '  -1.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.65'
$ws.Range("E26").Value = '  +1.66%  '
$ws.Range("E27").Value = '  +0.35%  '
$ws.Range("D28").Value = '1.958.37'
$ws.Range("E28").Value = '  -2.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.085'
$ws.Range("E29").Value = '  -4.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.89'
$ws.Range("E30").Value = '  -2.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.083'
$ws.Range("E31").Value = '  -7.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.658'
$ws.Range("E32").Value = '  +3.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08998'
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.527'
$ws.Range("E34").Value = '  -4.59%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.00'
$ws.Range("E35").Value = '  -5.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02298'
$ws.Range("E36").Value = '  -1.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2086'
$ws.Range("E37").Value = '  -1.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06009'
$ws.Range("E38").Value = '  -1.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6344'
$ws.Range("E39").Value = '  -2.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.948'
$ws.Range("E40").Value = '  -3.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.202'
$ws.Range("E41").Value = '  +1.20%  '
$ws.Range("B42").Value = 'WEMIXTOKEN'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.388'
$ws.Range("E42").Value = '  -2.43%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.771'
$ws.Range("E43").Value = '  -1.76%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.16'
$ws.Range("E44").Value = '  -2.67%  '
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.713'
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5886'
$ws.Range("E46").Value = '  -1.96%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '121.97'
$ws.Range("E47").Value = '  -2.35%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.952'
$ws.Range("E48").Value = '  -1.86%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.145'
$ws.Range("E49").Value = '  -0.66%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06857'
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.03'
$ws.Range("E51").Value = '  -3.70%  '
